$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the reporting period (row 8) -------------------------------
# Fecha de inicio del periodo que se informa: 2022-01-01 -> 2022-04-01
$ws.Range("B8").Value = 44652
# Fecha de término del periodo que se informa: 2022-03-31 -> 2022-06-30
$ws.Range("C8").Value = 44742
# Fecha de validación: 2022-04-08 -> 2022-07-11
$ws.Range("F8").Value = 44753
# Fecha de Actualización: 2022-04-08 -> 2022-07-11
$ws.Range("G8").Value = 44753

# --- Update the hyperlink's displayed text (D8) ------------------------
# Only the visible label changes (Enero-Marzo/f02 -> Abril-Junio/02_b);
# the underlying hyperlink target relationship is left untouched.
$ws.Range("D8").Value = "http://transparenciadocs.hidalgo.gob.mx/ENTIDADES/UPPachuca/dir1/2022/Abril-Junio/02_b/ESTRUCTURA%20ORGANICA%20%282%29.pdf"

# --- Update the selected cell shown when the workbook is reopened ------
$ws.Range("A9").Select()
